$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 605 - this pushes the existing rows 605-614
# down to 608-617 (unchanged), and leaves rows 605-607 blank (except for the
# date-formatted style on column D, which Excel carries down from row 605).
$ws.Rows("605:607").Insert()

# --- New row 605: Abate Fettel / Primera (O'Higgins) ---
$ws.Cells.Item(605, 1).Value = 7
$ws.Cells.Item(605, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(605, 3).Value = "Ñuble"
$ws.Cells.Item(605, 4).Value = 45239
$ws.Cells.Item(605, 5).Value = 16
$ws.Cells.Item(605, 6).Value = "Fruta"
$ws.Cells.Item(605, 7).Value = 100104
$ws.Cells.Item(605, 8).Value = "Frutos de pepita"
$ws.Cells.Item(605, 9).Value = 100104005
$ws.Cells.Item(605, 10).Value = "Pera"
$ws.Cells.Item(605, 11).Value = "Abate Fettel"
$ws.Cells.Item(605, 12).Value = "Primera"
$ws.Cells.Item(605, 13).Value = 80
$ws.Cells.Item(605, 14).Value = 15000
$ws.Cells.Item(605, 15).Value = 15000
$ws.Cells.Item(605, 16).Value = 15000
$ws.Cells.Item(605, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(605, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(605, 19).Value = 833
$ws.Cells.Item(605, 20).Value = 18

# --- New row 606: Packham's Triumph / Especial (Provincia de Curicó) ---
$ws.Cells.Item(606, 1).Value = 7
$ws.Cells.Item(606, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(606, 3).Value = "Ñuble"
$ws.Cells.Item(606, 4).Value = 45239
$ws.Cells.Item(606, 5).Value = 16
$ws.Cells.Item(606, 6).Value = "Fruta"
$ws.Cells.Item(606, 7).Value = 100104
$ws.Cells.Item(606, 8).Value = "Frutos de pepita"
$ws.Cells.Item(606, 9).Value = 100104005
$ws.Cells.Item(606, 10).Value = "Pera"
$ws.Cells.Item(606, 11).Value = "Packham's Triumph"
$ws.Cells.Item(606, 12).Value = "Especial"
$ws.Cells.Item(606, 13).Value = 100
$ws.Cells.Item(606, 14).Value = 17000
$ws.Cells.Item(606, 15).Value = 17000
$ws.Cells.Item(606, 16).Value = 17000
$ws.Cells.Item(606, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(606, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(606, 19).Value = 944
$ws.Cells.Item(606, 20).Value = 18

# --- New row 607: Packham's Triumph / Primera (Provincia de Curicó) ---
$ws.Cells.Item(607, 1).Value = 7
$ws.Cells.Item(607, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(607, 3).Value = "Ñuble"
$ws.Cells.Item(607, 4).Value = 45239
$ws.Cells.Item(607, 5).Value = 16
$ws.Cells.Item(607, 6).Value = "Fruta"
$ws.Cells.Item(607, 7).Value = 100104
$ws.Cells.Item(607, 8).Value = "Frutos de pepita"
$ws.Cells.Item(607, 9).Value = 100104005
$ws.Cells.Item(607, 10).Value = "Pera"
$ws.Cells.Item(607, 11).Value = "Packham's Triumph"
$ws.Cells.Item(607, 12).Value = "Primera"
$ws.Cells.Item(607, 13).Value = 100
$ws.Cells.Item(607, 14).Value = 15000
$ws.Cells.Item(607, 15).Value = 15000
$ws.Cells.Item(607, 16).Value = 15000
$ws.Cells.Item(607, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(607, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(607, 19).Value = 833
$ws.Cells.Item(607, 20).Value = 18
